$d = $word.ActiveDocument

# 1. "ithimbeni locwaningo" -> "ithimba locwaningo"
$d.Content.Find.Execute(
    "uthumele i-email ithimbeni locwaningo ku-",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "uthumele i-email ithimba locwaningo ku-",
    2)

# 2. "uzizwa ukhululekile" -> "uzizwe ukhululekile"
$d.Content.Find.Execute(
    "kuqinisekise ukuthi uzizwa ukhululekile",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "kuqinisekise ukuthi uzizwe ukhululekile",
    2)

# 3. "kanti iMenenja yocwaningo nguZamakhanya" -> "kanye neMenenja yocwaningo uZamakhanya"
$d.Content.Find.Execute(
    "kanti iMenenja yocwaningo nguZamakhanya",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "kanye neMenenja yocwaningo uZamakhanya",
    2)

# 4. "unemibuzo noma okukukhathazayo" -> "unemibuzo noma kukhona okukukhathazayo"
$d.Content.Find.Execute(
    "unemibuzo noma okukukhathazayo mayelana",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "unemibuzo noma kukhona okukukhathazayo mayelana",
    2)

# 5. "waqonda idokhumenti engenhla" -> "waqonda incwadi engenhla"
$d.Content.Find.Execute(
    "waqonda idokhumenti engenhla",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "waqonda incwadi engenhla",
    2)
